# Update daily metrics data: refresh firstName/lastName values and
# fix the swapped providerId values on rows 3 and 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# firstName column (B) - updated for all data rows
$ws.Range("B2").Value = "Kimberly"
$ws.Range("B3").Value = "Mark"
$ws.Range("B4").Value = "David"
$ws.Range("B5").Value = "Kayla"

# lastName column (C) - updated for all data rows
$ws.Range("C2").Value = "Floyd"
$ws.Range("C3").Value = "Gonzales"
$ws.Range("C4").Value = "Herring"
$ws.Range("C5").Value = "Mcdaniel"

# providerId column (A) - rows 3 and 4 swapped
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 1
